$xlShiftDown = -4121

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$table = $ws.ListObjects.Item("Table2")

# Insert a new row right after the current last table row (row 62), shifting
# cells down, which carries over the row's cell formatting/styles the same
# way Excel does when a table auto-expands.
$ws.Range("B62:F62").Copy()
$ws.Range("B63:F63").Insert($xlShiftDown)

# Fill in Post 53 values. Order mirrors how the source workbook's shared
# strings ended up appended (dev.to link, then title, then hashnode link).
$ws.Range("B63").Value = 53
$ws.Range("F63").Value = "https://dev.to/rahulmishra05/deadlock-prevention-operating-system-m04-p04-4khe"
$ws.Range("C63").Value = "Deadlock Prevention | Operating System - M04 P04"
$ws.Range("D63").Value = 44175
$ws.Range("E63").Value = "https://programmingport.hashnode.dev/deadlock-prevention-or-operating-system-m04-p04"

# Grow the table so its range/autofilter cover the new row.
$table.Resize($ws.Range("B10:F63"))

$ws.Range("E63").Select()
